$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.256.52"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.610.00"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "301.91"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "0.3777"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.3522"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "0.08069"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Value = "1.197"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "22.01"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("D14").Value = "6.348"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "7.239"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "0.00001208"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").Value = "1.590.26"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").Value = "93.97"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06900"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "17.14"
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").Value = "12.25"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "23.204.31"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "2.538"
$ws.Range("E25").Value = "  +3.76%  "
$ws.Range("D26").Value = "3.069"
$ws.Range("E26").Value = "  -5.82%  "
$ws.Range("D27").Value = "20.79"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "151.12"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").Value = "5.245"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").Value = "131.74"
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("D31").Value = "1.774.22"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "1.062"
$ws.Range("E32").Value = "  +11.46%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.421"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.097"
$ws.Range("E34").Value = "  -9.31%  "
$ws.Range("D35").Value = "11.32"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").Value = "0.02698"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").Value = "0.08669"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "0.06898"
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("D40").Value = "5.812"
$ws.Range("E40").Value = "  -4.74%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.311"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.6823"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("D43").Value = "11.94"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "15.22"
$ws.Range("E44").Value = "  -6.10%  "
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "0.6274"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").Value = "3.933"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "2.242"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("D49").Value = "0.07859"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "127.82"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").Value = "1.163"
$ws.Range("E51").Value = "  -3.16%  "
